# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (G) values for each data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 4
    8  = 1
    9  = 2
    10 = 1
    11 = 2
    12 = 0
    13 = 3
    14 = 1
    15 = 1
    16 = 3
    17 = 0
    18 = 0
    19 = 3
    20 = 2
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
